$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.038516640663147
$ws.Range("B1").Value = 1.794321417808533
$ws.Range("C1").Value = 6.852307796478271
$ws.Range("D1").Value = 1.624589920043945
$ws.Range("E1").Value = 0.9391082525253296
